# BF: typos and wrong places for some translations.
# Adds 23 new i18n strings (rows 320-340) to the "Feuil1" translation sheet,
# fixes a missing "Notes" annotation on row 311 (B311), and moves the
# selection/view to the newly appended area.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Step 1: create new shared strings in canonical order by writing their first-use cell
$ws.Range('A320').Value = 'The overtime you''ve requested has been rejected. Below, the details :'
$ws.Range('A321').Value = 'The overtime you''ve requested has been accepted. Below, the details :'
$ws.Range('A322').Value = '{Firstname} {Lastname} requests an overtime. Below, the details :'
$ws.Range('A324').Value = 'Once connected, you can change your password, as explained here.'
$ws.Range('A325').Value = 'The leave you''ve requested has been rejected. Below, the details :'
$ws.Range('A326').Value = 'The leave you''ve requested has been accepted. Below, the details :'
$ws.Range('A327').Value = '{Firstname} {Lastname} requests a leave. Below, the details :'
$ws.Range('A328').Value = 'From'
$ws.Range('A329').Value = 'To'
$ws.Range('A330').Value = 'Dear {Firstname} {Lastname},'
$ws.Range('A331').Value = 'If you didn''t perform this operation, please contact your administrator.'
$ws.Range('A332').Value = 'Welcome in LMS. If your are an employee, you could now :'
$ws.Range('A333').Value = 'See your leave balance.'
$ws.Range('A334').Value = 'See the list of the leave requests you have submitted.'
$ws.Range('A335').Value = 'Request a new leave.'
$ws.Range('A336').Value = 'If your are the line manager of other employee(s), you could now :'
$ws.Range('A337').Value = 'Validate leave requests submitted to you.'
$ws.Range('A338').Value = 'Validate overtime requests submitted to you.'
$ws.Range('A339').Value = 'Access forbidden'
$ws.Range('A340').Value = 'You are not allowed to perform this action.'
$ws.Range('B322').Value = 'don''t remove or replace {Firstname} and {Lastname}'
$ws.Range('A323').Value = 'Welcome to LMS {Firstname} {Lastname}. Please use these credentials to login to the system:'
$ws.Range('B311').Value = 'In the sense of overtime (extra hours)'

# Step 2: fill in duplicate cells that reuse already-created shared strings
$ws.Range('B323').Value = 'don''t remove or replace {Firstname} and {Lastname}'
$ws.Range('B330').Value = 'don''t remove or replace {Firstname} and {Lastname}'

# Step 3: update the view selection to match target state
[void]$ws.Range("A329").Select()
